$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.041080445316823
$ws.Range("C2").Value = 0.2251125788462787
$ws.Range("D2").Value = 0.01913225189683487
$ws.Range("F2").Value = 0.6445015020108613
$ws.Range("G2").Value = 0.002407840440158571
$ws.Range("I2").Value = 0.5678870897699433
$ws.Range("L2").Value = 0.3006213204614312
$ws.Range("N2").Value = 1.061035741188654
$ws.Range("O2").Value = 2.156353728437779
$ws.Range("B3").Value = 0.9386921363173428
$ws.Range("C3").Value = 0.2135994316831926
$ws.Range("D3").Value = 0.01779548436950051
$ws.Range("F3").Value = 0.6375846673627237
$ws.Range("G3").Value = 0.002410606811406165
$ws.Range("I3").Value = 0.5723319223654393
$ws.Range("L3").Value = 0.2891276646306125
$ws.Range("N3").Value = 1.068791596941715
$ws.Range("O3").Value = 2.147307025768157
$ws.Range("B4").Value = 0.875907542782727
$ws.Range("C4").Value = 0.2064758940172027
$ws.Range("D4").Value = 0.01696855492318861
$ws.Range("F4").Value = 0.6338077703389757
$ws.Range("G4").Value = 0.002412397351715247
$ws.Range("I4").Value = 0.5754438552258314
$ws.Range("L4").Value = 0.2822435626249558
$ws.Range("N4").Value = 1.07399297515007
$ws.Range("O4").Value = 2.143310348887184
$ws.Range("B5").Value = 0.8503445157002716
$ws.Range("C5").Value = 0.2035594932028033
$ws.Range("D5").Value = 0.01663004875400276
$ws.Range("F5").Value = 0.6323867316214375
$ws.Range("G5").Value = 0.002413150209214658
$ws.Range("I5").Value = 0.5768081263585572
$ws.Range("L5").Value = 0.279481786882684
$ws.Range("N5").Value = 1.07622313339629
$ws.Range("O5").Value = 2.142072949049862
$ws.Range("B6").Value = 0.8461011805705994
$ws.Range("C6").Value = 0.2030744170031511
$ws.Range("D6").Value = 0.01657374851601645
$ws.Range("F6").Value = 0.6321578971685327
$ws.Range("G6").Value = 0.00241327662372618
$ws.Range("I6").Value = 0.577040465183984
$ws.Range("L6").Value = 0.2790258276516084
$ws.Range("N6").Value = 1.076600129746254
$ws.Range("O6").Value = 2.141891097325328
$ws.Range("B7").Value = 0.8755626989329812
$ws.Range("C7").Value = 0.2064366168396816
$ws.Range("D7").Value = 0.01696399586123931
$ws.Range("F7").Value = 0.6337881277892663
$ws.Range("G7").Value = 0.002412407410950554
$ws.Range("I7").Value = 0.5754618651863588
$ws.Range("L7").Value = 0.2822061399617155
$ws.Range("N7").Value = 1.074022604077769
$ws.Range("O7").Value = 2.143292077278886
$ws.Range("B8").Value = 1.005760768505468
$ws.Range("C8").Value = 0.2211542465107641
$ws.Range("D8").Value = 0.01867262251724355
$ws.Range("F8").Value = 0.642018935821234
$ws.Range("G8").Value = 0.002408775236865668
$ws.Range("I8").Value = 0.5693401467490133
$ws.Range("L8").Value = 0.2966224095376049
$ws.Range("N8").Value = 1.063618852779314
$ws.Range("O8").Value = 2.152910743598056
$ws.Range("B9").Value = 1.261677102275371
$ws.Range("C9").Value = 0.249576933680288
$ws.Range("D9").Value = 0.02197369494924573
$ws.Range("F9").Value = 0.661896804175683
$ws.Range("G9").Value = 0.002402379194560689
$ws.Range("I9").Value = 0.5603785055491031
$ws.Range("L9").Value = 0.3262659472450764
$ws.Range("N9").Value = 1.046698293077576
$ws.Range("O9").Value = 2.184162362785088
$ws.Range("B10").Value = 1.450009435968582
$ws.Range("C10").Value = 0.270184538136931
$ws.Range("D10").Value = 0.02436799928319289
$ws.Range("F10").Value = 0.6787925108250903
$ws.Range("G10").Value = 0.002398118580093426
$ws.Range("I10").Value = 0.5556579149273801
$ws.Range("L10").Value = 0.3488852465508927
$ws.Range("N10").Value = 1.03638396106129
$ws.Range("O10").Value = 2.214719883207522
$ws.Range("B11").Value = 1.535743821801589
$ws.Range("C11").Value = 0.2794984459568184
$ws.Range("D11").Value = 0.02545034514363209
$ws.Range("F11").Value = 0.6869793362357655
$ws.Range("G11").Value = 0.002396274604306952
$ws.Range("I11").Value = 0.553916884569233
$ws.Range("L11").Value = 0.3593585967682884
$ws.Range("N11").Value = 1.032150497751502
$ws.Range("O11").Value = 2.230280745571065
$ws.Range("B12").Value = 1.568216693877446
$ws.Range("C12").Value = 0.2830165058072964
$ws.Range("D12").Value = 0.02585920039842904
$ws.Range("F12").Value = 0.6901516820536955
$ws.Range("G12").Value = 0.00239558981332966
$ws.Range("I12").Value = 0.5533161848748023
$ws.Range("L12").Value = 0.363351005192655
$ws.Range("N12").Value = 1.030613271502446
$ws.Range("O12").Value = 2.236412623391288
$ws.Range("B13").Value = 1.561222789553312
$ws.Range("C13").Value = 0.2822592283678489
$ws.Range("D13").Value = 0.0257711911913816
$ws.Range("F13").Value = 0.6894652484070889
$ws.Range("G13").Value = 0.002395736696766184
$ws.Range("I13").Value = 0.5534429480418055
$ws.Range("L13").Value = 0.362489996085003
$ws.Range("N13").Value = 1.030941410729625
$ws.Range("O13").Value = 2.235081360619404
$ws.Range("B14").Value = 1.538415251313893
$ws.Range("C14").Value = 0.2797880586780082
$ws.Range("D14").Value = 0.02548400219579605
$ws.Range("F14").Value = 0.6872388797543749
$ws.Range("G14").Value = 0.002396217996183174
$ws.Range("I14").Value = 0.5538662895273205
$ws.Range("L14").Value = 0.3596865261907141
$ws.Range("N14").Value = 1.032022708861838
$ws.Range("O14").Value = 2.230780419169889
$ws.Range("B15").Value = 1.524445853726945
$ws.Range("C15").Value = 0.2782732289375076
$ws.Range("D15").Value = 0.02530795900656102
$ws.Range("F15").Value = 0.685884568002578
$ws.Range("G15").Value = 0.002396514560263885
$ws.Range("I15").Value = 0.5541332331670361
$ws.Range("L15").Value = 0.3579727539191282
$ws.Range("N15").Value = 1.032693615229128
$ws.Range("O15").Value = 2.228177152159333
$ws.Range("B16").Value = 1.44440757899099
$ws.Range("C16").Value = 0.2695746158011616
$ws.Range("D16").Value = 0.02429712594998534
$ws.Range("F16").Value = 0.678267574126167
$ws.Range("G16").Value = 0.002398240978223786
$ws.Range("I16").Value = 0.5557798859905887
$ws.Range("L16").Value = 0.3482044819117789
$ws.Range("N16").Value = 1.036669851350062
$ws.Range("O16").Value = 2.213736400727868
$ws.Range("B17").Value = 1.395321147362495
$ws.Range("C17").Value = 0.2642226351044599
$ws.Range("D17").Value = 0.02367524576442293
$ws.Range("F17").Value = 0.6737231879502588
$ws.Range("G17").Value = 0.002399324161311298
$ws.Range("I17").Value = 0.5568942488813207
$ws.Range("L17").Value = 0.3422589893463339
$ws.Range("N17").Value = 1.039226560219795
$ws.Range("O17").Value = 2.205303083311065
$ws.Range("B18").Value = 1.367093776680974
$ws.Range("C18").Value = 0.2611386259014239
$ws.Range("D18").Value = 0.02331691472612363
$ws.Range("F18").Value = 0.6711565141078779
$ws.Range("G18").Value = 0.002399956049687947
$ws.Range("I18").Value = 0.5575734432018891
$ws.Range("L18").Value = 0.3388565940871757
$ws.Range("N18").Value = 1.040740277668476
$ws.Range("O18").Value = 2.200608682456817
$ws.Range("B19").Value = 1.357537528758996
$ws.Range("C19").Value = 0.2600934620082569
$ws.Range("D19").Value = 0.02319548042239461
$ws.Range("F19").Value = 0.6702955737748368
$ws.Range("G19").Value = 0.002400171521870341
$ws.Range("I19").Value = 0.5578099702285755
$ws.Range("L19").Value = 0.3377075731344519
$ws.Range("N19").Value = 1.041260211770712
$ws.Range("O19").Value = 2.199046051718057
$ws.Range("B20").Value = 1.400545888748752
$ws.Range("C20").Value = 0.2647929530923818
$ws.Range("D20").Value = 0.02374151258131718
$ws.Range("F20").Value = 0.6742020660605021
$ws.Range("G20").Value = 0.00239920793710147
$ws.Range("I20").Value = 0.556771664013219
$ws.Range("L20").Value = 0.342890107523445
$ws.Range("N20").Value = 1.038949927028483
$ws.Range("O20").Value = 2.206184651792341
$ws.Range("B21").Value = 1.545114198112401
$ws.Range("C21").Value = 0.2805141446676771
$ws.Range("D21").Value = 0.02556838403985751
$ws.Range("F21").Value = 0.6878908582906575
$ws.Range("G21").Value = 0.002396076261421883
$ws.Range("I21").Value = 0.5537403525233699
$ws.Range("L21").Value = 0.3605092576890172
$ws.Range("N21").Value = 1.031703317396797
$ws.Range("O21").Value = 2.232037210188935
$ws.Range("B22").Value = 1.639638683284431
$ws.Range("C22").Value = 0.2907367961933289
$ws.Range("D22").Value = 0.02675647802089287
$ws.Range("F22").Value = 0.6972580199746403
$ws.Range("G22").Value = 0.002394108088276028
$ws.Range("I22").Value = 0.552100797226295
$ws.Range("L22").Value = 0.3721781229082097
$ws.Range("N22").Value = 1.027351303053564
$ws.Range("O22").Value = 2.250328549098612
$ws.Range("B23").Value = 1.589186006659872
$ws.Range("C23").Value = 0.2852856042293581
$ws.Range("D23").Value = 0.02612291470812522
$ws.Range("F23").Value = 0.6922200433878203
$ws.Range("G23").Value = 0.002395151372122454
$ws.Range("I23").Value = 0.5529445537247142
$ws.Range("L23").Value = 0.3659361769092442
$ws.Range("N23").Value = 1.029638927953556
$ws.Range("O23").Value = 2.240438259969949
$ws.Range("B24").Value = 1.39818380461935
$ws.Range("C24").Value = 0.2645351343822426
$ws.Range("D24").Value = 0.02371155585541374
$ws.Range("F24").Value = 0.673985422132958
$ws.Range("G24").Value = 0.002399260453750511
$ws.Range("I24").Value = 0.5568269646166897
$ws.Range("L24").Value = 0.3426047299428632
$ws.Range("N24").Value = 1.039074856355732
$ws.Range("O24").Value = 2.205785614943153
$ws.Range("B25").Value = 1.192386510976576
$ws.Range("C25").Value = 0.2419355385580673
$ws.Range("D25").Value = 0.02108605302841937
$ws.Range("F25").Value = 0.6561178550236164
$ws.Range("G25").Value = 0.002404032160699262
$ws.Range("I25").Value = 0.5624762192392794
$ws.Range("L25").Value = 0.3180992538901535
$ws.Range("N25").Value = 1.030941410729625
$ws.Range("O25").Value = 2.235081360619404
